# Auto-generated edit script applying the diff to Midgardsormr_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 17862172
$ws.Range("I51").Value = 41669830
$ws.Range("J51").Value = 6424.75
$ws.Range("K51").Value = 41669830
$ws.Range("L51").Value = 6424.75
$ws.Range("M51").Value = -41669346
$ws.Range("N51").Value = -7392.75
$ws.Range("H62").Value = 8051.8125
$ws.Range("I62").Value = 8510.429
$ws.Range("K62").Value = 8510.429
$ws.Range("M62").Value = -7886.429
$ws.Range("H65").Value = 8051.8125
$ws.Range("I65").Value = 8510.429
$ws.Range("K65").Value = 42552.145
$ws.Range("M65").Value = -39432.145
$ws.Range("H74").Value = 2320.1428
$ws.Range("I74").Value = 2320.1428
$ws.Range("K74").Value = 2320.1428
$ws.Range("M74").Value = -1384.1428
$ws.Range("H77").Value = 2320.1428
$ws.Range("I77").Value = 2320.1428
$ws.Range("K77").Value = 11600.714
$ws.Range("M77").Value = -6920.714
$ws.Range("H100").Value = 37826.234
$ws.Range("I100").Value = 61203.2
$ws.Range("K100").Value = 61203.2
$ws.Range("M100").Value = -60662.2
$ws.Range("H111").Value = 899.6
$ws.Range("J111").Value = 1209.4286
$ws.Range("L111").Value = 3628.2858
$ws.Range("N111").Value = -9762.2858
$ws.Range("H116").Value = 8803.909
$ws.Range("I116").Value = 4830.6665
$ws.Range("K116").Value = 4830.6665
$ws.Range("M116").Value = -1388.6665
$ws.Range("H132").Value = 2659738.5
$ws.Range("I132").Value = 2659738.5
$ws.Range("K132").Value = 7979215.5
$ws.Range("M132").Value = -7976685.5
$ws.Range("H134").Value = 105666.664
$ws.Range("J134").Value = 105666.664
$ws.Range("L134").Value = 105666.664
$ws.Range("N134").Value = -115806.664
$ws.Range("H137").Value = 15522.125
$ws.Range("J137").Value = 3464.5833
$ws.Range("L137").Value = 10393.7499
$ws.Range("N137").Value = -15493.7499
$ws.Range("H139").Value = 67000
$ws.Range("J139").Value = 67000
$ws.Range("L139").Value = 67000
$ws.Range("N139").Value = -77280
$ws.Range("H140").Value = 67000
$ws.Range("J140").Value = 67000
$ws.Range("L140").Value = 67000
$ws.Range("N140").Value = -77360
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2844.8125
$ws.Range("I45").Value = 2050.5518
$ws.Range("J45").Value = 4057.1052
$ws.Range("K45").Value = 2050.5518
$ws.Range("L45").Value = 4057.1052
$ws.Range("M45").Value = -1673.5518
$ws.Range("N45").Value = -4811.1052
$ws.Range("H55").Value = 49998.5
$ws.Range("J55").Value = 49998.5
$ws.Range("L55").Value = 49998.5
$ws.Range("N55").Value = -50628.5
$ws.Range("H88").Value = 9905.571
$ws.Range("I88").Value = 3467.8
$ws.Range("J88").Value = 26000
$ws.Range("K88").Value = 3467.8
$ws.Range("L88").Value = 26000
$ws.Range("M88").Value = -3061.8
$ws.Range("N88").Value = -26812
$ws.Range("H91").Value = 9905.571
$ws.Range("I91").Value = 3467.8
$ws.Range("J91").Value = 26000
$ws.Range("K91").Value = 3467.8
$ws.Range("L91").Value = 26000
$ws.Range("M91").Value = -2063.8
$ws.Range("N91").Value = -28808
$ws.Range("H122").Value = 2619.457
$ws.Range("I122").Value = 2485.9033
$ws.Range("K122").Value = 7457.7099
$ws.Range("M122").Value = -5007.7099
$ws.Range("H134").Value = 79998.664
$ws.Range("J134").Value = 79998.664
$ws.Range("L134").Value = 79998.664
$ws.Range("N134").Value = -90138.664
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2239.1333
$ws.Range("I105").Value = 1259.8334
$ws.Range("K105").Value = 1259.8334
$ws.Range("M105").Value = 487.1666
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1829.2142
$ws.Range("I16").Value = 1361.9
$ws.Range("K16").Value = 1361.9
$ws.Range("M16").Value = -1074.9
$ws.Range("H31").Value = 2441338.2
$ws.Range("I31").Value = 3228078.5
$ws.Range("J31").Value = 2443.7
$ws.Range("K31").Value = 3228078.5
$ws.Range("L31").Value = 2443.7
$ws.Range("M31").Value = -3227783.5
$ws.Range("N31").Value = -3033.7
$ws.Range("H34").Value = 2441338.2
$ws.Range("I34").Value = 3228078.5
$ws.Range("J34").Value = 2443.7
$ws.Range("K34").Value = 3228078.5
$ws.Range("L34").Value = 2443.7
$ws.Range("M34").Value = -3227876.5
$ws.Range("N34").Value = -2847.7
$ws.Range("H107").Value = 564.225
$ws.Range("I107").Value = 452.36
$ws.Range("K107").Value = 452.36
$ws.Range("M107").Value = 1467.64
$ws.Range("H113").Value = 1829.2142
$ws.Range("I113").Value = 1361.9
$ws.Range("K113").Value = 1361.9
$ws.Range("M113").Value = 808.0999999999999
$ws.Range("H138").Value = 115451.586
$ws.Range("J138").Value = 115451.586
$ws.Range("L138").Value = 115451.586
$ws.Range("N138").Value = -125731.586
$ws.Range("H140").Value = 104410.07
$ws.Range("J140").Value = 104410.07
$ws.Range("L140").Value = 104410.07
$ws.Range("N140").Value = -114770.07
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 2489.2222
$ws.Range("J41").Value = 3667.1667
$ws.Range("L41").Value = 11001.5001
$ws.Range("N41").Value = -11677.5001
$ws.Range("H56").Value = 7776
$ws.Range("I56").Value = 7776
$ws.Range("K56").Value = 7776
$ws.Range("M56").Value = -7246
$ws.Range("H81").Value = 4916.6665
$ws.Range("I81").Value = 3250
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 9750
$ws.Range("L81").Value = 15000
$ws.Range("M81").Value = -8627
$ws.Range("N81").Value = -17246
$ws.Range("H84").Value = 4916.6665
$ws.Range("I84").Value = 3250
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 29250
$ws.Range("L84").Value = 45000
$ws.Range("M84").Value = -23634
$ws.Range("N84").Value = -56232
$ws.Range("H132").Value = 2026.2354
$ws.Range("I132").Value = 2094.875
$ws.Range("J132").Value = 1965.2222
$ws.Range("K132").Value = 18853.875
$ws.Range("L132").Value = 17686.9998
$ws.Range("M132").Value = -16323.875
$ws.Range("N132").Value = -22746.9998
$ws.Range("H136").Value = 5460.643
$ws.Range("I136").Value = 2406.125
$ws.Range("J136").Value = 9533.333000000001
$ws.Range("K136").Value = 7218.375
$ws.Range("L136").Value = 28599.999
$ws.Range("M136").Value = -2118.375
$ws.Range("N136").Value = -38799.999
$ws.Range("H138").Value = 12083.786
$ws.Range("I138").Value = 13647.5
$ws.Range("K138").Value = 40942.5
$ws.Range("M138").Value = -35802.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").Value = ""
$ws.Range("H122").Value = 3172.8823
$ws.Range("I122").Value = 3172.8823
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9518.6469
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7068.6469
$ws.Range("N122").Value = ""
$ws.Range("H139").Value = 104995.336
$ws.Range("J139").Value = 104995.336
$ws.Range("L139").Value = 104995.336
$ws.Range("N139").Value = -115275.336
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 57243676
$ws.Range("I12").Value = 68691610
$ws.Range("K12").Value = 68691610
$ws.Range("M12").Value = -68691440
$ws.Range("H22").Value = 1827
$ws.Range("J22").Value = 1955.7778
$ws.Range("L22").Value = 1955.7778
$ws.Range("N22").Value = -2545.7778
$ws.Range("H27").Value = 1827
$ws.Range("J27").Value = 1955.7778
$ws.Range("L27").Value = 1955.7778
$ws.Range("N27").Value = -2169.7778
$ws.Range("H46").Value = 3399.8965
$ws.Range("I46").Value = 785.2308
$ws.Range("J46").Value = 5524.3125
$ws.Range("K46").Value = 785.2308
$ws.Range("L46").Value = 5524.3125
$ws.Range("M46").Value = -597.2308
$ws.Range("N46").Value = -5900.3125
$ws.Range("H55").Value = 984.4074000000001
$ws.Range("I55").Value = 364.72726
$ws.Range("J55").Value = 1410.4375
$ws.Range("K55").Value = 364.72726
$ws.Range("L55").Value = 1410.4375
$ws.Range("M55").Value = -191.72726
$ws.Range("N55").Value = -1756.4375
$ws.Range("H122").Value = 8985.290000000001
$ws.Range("I122").Value = 9939.526
$ws.Range("K122").Value = 29818.578
$ws.Range("M122").Value = -27368.578
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3350.6316
$ws.Range("I81").Value = 3655.0715
$ws.Range("K81").Value = 7310.143
$ws.Range("M81").Value = -6249.143
$ws.Range("H84").Value = 3350.6316
$ws.Range("I84").Value = 3655.0715
$ws.Range("K84").Value = 36550.715
$ws.Range("M84").Value = -31246.715
$ws.Range("H132").Value = 3786.5938
$ws.Range("I132").Value = 4066.682
$ws.Range("J132").Value = 3170.4
$ws.Range("K132").Value = 12200.046
$ws.Range("L132").Value = 9511.200000000001
$ws.Range("M132").Value = -9670.045999999998
$ws.Range("N132").Value = -14571.2
$ws.Range("H136").Value = 42411.438
$ws.Range("I136").Value = 45045.535
$ws.Range("K136").Value = 135136.605
$ws.Range("M136").Value = -132586.605
